$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '75.898.94'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.882.25'
$ws.Range("E3").Value = '  +5.80%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '194.89'
$ws.Range("E5").Value = '  +3.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '596.83'
$ws.Range("E6").Value = '  +0.76%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.192'
$ws.Range("E9").Value = '  -2.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2.882.46'
$ws.Range("E10").Value = '  +5.83%  '
$ws.Range("E11").Value = '  +9.12%  '
$ws.Range("E12").Value = '  -1.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.89'
$ws.Range("E13").Value = '  +1.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.351.79'
$ws.Range("E14").Value = '  +3.85%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '75.772.94'
$ws.Range("E15").Value = '  +0.26%  '
$ws.Range("E16").Value = '  -1.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.24'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.873.12'
$ws.Range("E18").Value = '  +5.21%  '
$ws.Range("E19").Value = '  -6.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.54'
$ws.Range("E20").Value = '  +2.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '376.59'
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.29'
$ws.Range("E22").Value = '  -1.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.14'
$ws.Range("E23").Value = '  +0.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.47'
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.026.88'
$ws.Range("E26").Value = '  +5.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.22'
$ws.Range("E27").Value = '  -0.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.73'
$ws.Range("E28").Value = '  +0.70%  '
$ws.Range("E29").Value = '  +6.86%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.36%  '
$ws.Range("E31").Value = '  -1.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '506.18'
$ws.Range("E32").Value = '  -3.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.76'
$ws.Range("E33").Value = '  -1.99%  '
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '20.18'
$ws.Range("E36").Value = '  +2.89%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.20'
$ws.Range("E37").Value = '  +1.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.70'
$ws.Range("E38").Value = '  +1.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.114'
$ws.Range("E39").Value = '  -5.46%  '
$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '181.55'
$ws.Range("E41").Value = '  +4.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.344'
$ws.Range("E42").Value = '  +2.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.98'
$ws.Range("E43").Value = '  -2.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.68'
$ws.Range("E44").Value = '  -3.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0915'
$ws.Range("E45").Value = '  +6.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.22'
$ws.Range("E46").Value = '  -1.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.41'
$ws.Range("E47").Value = '  +3.12%  '
$ws.Range("E48").Value = '  -3.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.578'
$ws.Range("E49").Value = '  +4.61%  '
$ws.Range("E50").Value = '  +11.13%  '
$ws.Range("E51").Value = '  +0.37%  '
